$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the header row (Username/Password) with a new data row (Admin/Paswword123)
$ws.Range("A1").Value = "Admin"
$ws.Range("B1").Value = "Paswword123"

# Move the selection to B1
$ws.Range("B1").Select()
